# Apply the edit described by the diff:
#  - Remove column C ("description"); column D ("reason") shifts left to become column C.
#  - Update the A (id) and B (score) values for rows 2-5.
#  - Replace the C (reason) text for rows 2-5 with updated wording.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old "description" column (C). This shifts the old "reason"
# column (D) left into column C and updates the used range / dimension.
$ws.Columns.Item(3).Delete()

function Set-TextValue($cell, $text) {
    # Force the cell to be stored as text (matching the source file, where
    # the "id" column holds numeric-looking values as strings) without
    # leaving a numeric style behind: the leading apostrophe makes Excel
    # treat the input as text, then resetting the style to "Normal"
    # clears the quote-prefix formatting flag that the apostrophe implies.
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Cells.Item(2, 1) "3"
$ws.Cells.Item(2, 2).Value = 99
$ws.Cells.Item(2, 3).Value = "The job of an NLP Engineer requires strong skills in NLP and Python, which you possess. Your experience in developing bias/toxicity detection pipelines and implementing language models aligns well with the job description. The high score indicates a strong match between your project experiences and the job requirements."

# Row 3
Set-TextValue $ws.Cells.Item(3, 1) "1"
$ws.Cells.Item(3, 2).Value = 79
$ws.Cells.Item(3, 3).Value = "The job as a SDE Intern with a score of 79.0 is suitable for you because it requires skills in MongoDB, ReactJS, JavaScript, Web Development, and NodeJS. Your experience in developing web applications using ReactJS, MongoDB, and NodeJS aligns well with the job requirements. Additionally, your skillset in JavaScript and web development will be valuable for code maintenance, scalability, and feature development."

# Row 4
Set-TextValue $ws.Cells.Item(4, 1) "4"
$ws.Cells.Item(4, 2).Value = 36.5
$ws.Cells.Item(4, 3).Value = "The job description of Application Development has a low score of 36.5. It is not suitable for you because it requires skills in Flutter, Dart, Android/iOS, Firestore, and Firebase Authentication which are not mentioned in your projects or resume."

# Row 5
Set-TextValue $ws.Cells.Item(5, 1) "2"
$ws.Cells.Item(5, 2).Value = 70
$ws.Cells.Item(5, 3).Value = "The job as a Frontend Engineer Intern requires proficiency in ReactJS, JavaScript, CSS, and NextJS, which aligns with your skills. The project 'LLMGuard' demonstrates your experience with ReactJS and JavaScript, making you suitable for this role. However, the unpaid evaluative internship and short duration may make it less desirable compared to other options."
